$wb = $excel.ActiveWorkbook

# --- Update time_taken (col F) timestamps on the "data" sheet ---
$ws1 = $wb.Worksheets.Item("data")

$ws1.Cells.Item(2, 6).Value = "2021-10-05 14:20:45.215835"
$ws1.Cells.Item(3, 6).Value = "2021-10-05 14:20:45.215843"
$ws1.Cells.Item(4, 6).Value = "2021-10-05 14:20:45.215846"
$ws1.Cells.Item(5, 6).Value = "2021-10-05 14:20:45.215848"
$ws1.Cells.Item(6, 6).Value = "2021-10-05 14:20:45.215851"
$ws1.Cells.Item(7, 6).Value = "2021-10-05 14:20:45.215854"
$ws1.Cells.Item(8, 6).Value = "2021-10-05 14:20:45.215856"
$ws1.Cells.Item(9, 6).Value = "2021-10-05 14:20:45.215859"
$ws1.Cells.Item(10, 6).Value = "2021-10-05 14:20:45.215862"
$ws1.Cells.Item(11, 6).Value = "2021-10-05 14:20:45.215864"
$ws1.Cells.Item(12, 6).Value = "2021-10-05 14:20:45.215867"
$ws1.Cells.Item(13, 6).Value = "2021-10-05 14:20:45.215869"
$ws1.Cells.Item(14, 6).Value = "2021-10-05 14:20:45.215872"
$ws1.Cells.Item(15, 6).Value = "2021-10-05 14:20:45.215874"
$ws1.Cells.Item(16, 6).Value = "2021-10-05 14:20:45.215876"
$ws1.Cells.Item(17, 6).Value = "2021-10-05 14:20:45.215879"
$ws1.Cells.Item(18, 6).Value = "2021-10-05 14:20:45.215882"
$ws1.Cells.Item(19, 6).Value = "2021-10-05 14:20:45.215884"
$ws1.Cells.Item(20, 6).Value = "2021-10-05 14:20:45.215887"
$ws1.Cells.Item(21, 6).Value = "2021-10-05 14:20:45.215890"
$ws1.Cells.Item(22, 6).Value = "2021-10-05 14:20:45.215892"
$ws1.Cells.Item(23, 6).Value = "2021-10-05 14:20:45.215895"
$ws1.Cells.Item(24, 6).Value = "2021-10-05 14:20:45.215897"
$ws1.Cells.Item(25, 6).Value = "2021-10-05 14:20:45.215900"
$ws1.Cells.Item(26, 6).Value = "2021-10-05 14:20:45.215903"
$ws1.Cells.Item(27, 6).Value = "2021-10-05 14:20:45.215905"
$ws1.Cells.Item(28, 6).Value = "2021-10-05 14:20:45.215908"
$ws1.Cells.Item(29, 6).Value = "2021-10-05 14:20:45.215910"
$ws1.Cells.Item(30, 6).Value = "2021-10-05 14:20:45.215913"
$ws1.Cells.Item(31, 6).Value = "2021-10-05 14:20:45.215915"
$ws1.Cells.Item(32, 6).Value = "2021-10-05 14:20:45.215917"
$ws1.Cells.Item(33, 6).Value = "2021-10-05 14:20:45.215920"
$ws1.Cells.Item(34, 6).Value = "2021-10-05 14:20:45.215923"
$ws1.Cells.Item(35, 6).Value = "2021-10-05 14:20:45.215926"
$ws1.Cells.Item(36, 6).Value = "2021-10-05 14:20:45.215928"
$ws1.Cells.Item(37, 6).Value = "2021-10-05 14:20:45.215930"
$ws1.Cells.Item(38, 6).Value = "2021-10-05 14:20:45.215933"
$ws1.Cells.Item(39, 6).Value = "2021-10-05 14:20:45.215935"
$ws1.Cells.Item(40, 6).Value = "2021-10-05 14:20:45.215938"
$ws1.Cells.Item(41, 6).Value = "2021-10-05 14:20:45.215940"
$ws1.Cells.Item(42, 6).Value = "2021-10-05 14:20:45.215943"
$ws1.Cells.Item(43, 6).Value = "2021-10-05 14:20:45.215946"
$ws1.Cells.Item(44, 6).Value = "2021-10-05 14:20:45.215948"
$ws1.Cells.Item(45, 6).Value = "2021-10-05 14:20:45.215950"
$ws1.Cells.Item(46, 6).Value = "2021-10-05 14:20:45.215953"
$ws1.Cells.Item(47, 6).Value = "2021-10-05 14:20:45.215955"
$ws1.Cells.Item(48, 6).Value = "2021-10-05 14:20:45.215958"
$ws1.Cells.Item(49, 6).Value = "2021-10-05 14:20:45.215960"
$ws1.Cells.Item(50, 6).Value = "2021-10-05 14:20:45.215963"
$ws1.Cells.Item(51, 6).Value = "2021-10-05 14:20:45.215965"
$ws1.Cells.Item(52, 6).Value = "2021-10-05 14:20:45.215967"
$ws1.Cells.Item(53, 6).Value = "2021-10-05 14:20:45.215970"
$ws1.Cells.Item(54, 6).Value = "2021-10-05 14:20:45.215973"
$ws1.Cells.Item(55, 6).Value = "2021-10-05 14:20:45.215975"
$ws1.Cells.Item(56, 6).Value = "2021-10-05 14:20:45.215978"
$ws1.Cells.Item(57, 6).Value = "2021-10-05 14:20:45.215980"
$ws1.Cells.Item(58, 6).Value = "2021-10-05 14:20:45.215982"
$ws1.Cells.Item(59, 6).Value = "2021-10-05 14:20:45.215985"
$ws1.Cells.Item(60, 6).Value = "2021-10-05 14:20:45.215988"
$ws1.Cells.Item(61, 6).Value = "2021-10-05 14:20:45.215990"
$ws1.Cells.Item(62, 6).Value = "2021-10-05 14:20:45.215992"
$ws1.Cells.Item(63, 6).Value = "2021-10-05 14:20:45.215995"
$ws1.Cells.Item(64, 6).Value = "2021-10-05 14:20:45.215997"
$ws1.Cells.Item(65, 6).Value = "2021-10-05 14:20:45.216000"
$ws1.Cells.Item(66, 6).Value = "2021-10-05 14:20:45.216003"
$ws1.Cells.Item(67, 6).Value = "2021-10-05 14:20:45.216006"
$ws1.Cells.Item(68, 6).Value = "2021-10-05 14:20:45.216008"
$ws1.Cells.Item(69, 6).Value = "2021-10-05 14:20:45.216011"
$ws1.Cells.Item(70, 6).Value = "2021-10-05 14:20:45.216013"
$ws1.Cells.Item(71, 6).Value = "2021-10-05 14:20:45.216016"
$ws1.Cells.Item(72, 6).Value = "2021-10-05 14:20:45.216018"
$ws1.Cells.Item(73, 6).Value = "2021-10-05 14:20:45.216021"
$ws1.Cells.Item(74, 6).Value = "2021-10-05 14:20:45.216023"
$ws1.Cells.Item(75, 6).Value = "2021-10-05 14:20:45.216026"
$ws1.Cells.Item(76, 6).Value = "2021-10-05 14:20:45.216028"
$ws1.Cells.Item(77, 6).Value = "2021-10-05 14:20:45.216031"
$ws1.Cells.Item(78, 6).Value = "2021-10-05 14:20:45.216035"
$ws1.Cells.Item(79, 6).Value = "2021-10-05 14:20:45.216038"
$ws1.Cells.Item(80, 6).Value = "2021-10-05 14:20:45.216041"
$ws1.Cells.Item(81, 6).Value = "2021-10-05 14:20:45.216043"
$ws1.Cells.Item(82, 6).Value = "2021-10-05 14:20:45.216046"
$ws1.Cells.Item(83, 6).Value = "2021-10-05 14:20:45.216048"
$ws1.Cells.Item(84, 6).Value = "2021-10-05 14:20:45.216051"
$ws1.Cells.Item(85, 6).Value = "2021-10-05 14:20:45.216053"
$ws1.Cells.Item(86, 6).Value = "2021-10-05 14:20:45.216056"
$ws1.Cells.Item(87, 6).Value = "2021-10-05 14:20:45.216058"
$ws1.Cells.Item(88, 6).Value = "2021-10-05 14:20:45.216060"
$ws1.Cells.Item(89, 6).Value = "2021-10-05 14:20:45.216063"
$ws1.Cells.Item(90, 6).Value = "2021-10-05 14:20:45.216065"
$ws1.Cells.Item(91, 6).Value = "2021-10-05 14:20:45.216068"
$ws1.Cells.Item(92, 6).Value = "2021-10-05 14:20:45.216070"
$ws1.Cells.Item(93, 6).Value = "2021-10-05 14:20:45.216073"
$ws1.Cells.Item(94, 6).Value = "2021-10-05 14:20:45.216076"
$ws1.Cells.Item(95, 6).Value = "2021-10-05 14:20:45.216079"
$ws1.Cells.Item(96, 6).Value = "2021-10-05 14:20:45.216082"
$ws1.Cells.Item(97, 6).Value = "2021-10-05 14:20:45.216084"
$ws1.Cells.Item(98, 6).Value = "2021-10-05 14:20:45.216087"
$ws1.Cells.Item(99, 6).Value = "2021-10-05 14:20:45.216089"
$ws1.Cells.Item(100, 6).Value = "2021-10-05 14:20:45.216092"
$ws1.Cells.Item(101, 6).Value = "2021-10-05 14:20:45.216094"
$ws1.Cells.Item(102, 6).Value = "2021-10-05 14:20:45.216097"
$ws1.Cells.Item(103, 6).Value = "2021-10-05 14:20:45.216099"
$ws1.Cells.Item(104, 6).Value = "2021-10-05 14:20:45.216102"
$ws1.Cells.Item(105, 6).Value = "2021-10-05 14:20:45.216104"
$ws1.Cells.Item(106, 6).Value = "2021-10-05 14:20:45.216107"
$ws1.Cells.Item(107, 6).Value = "2021-10-05 14:20:45.216109"
$ws1.Cells.Item(108, 6).Value = "2021-10-05 14:20:45.216111"
$ws1.Cells.Item(109, 6).Value = "2021-10-05 14:20:45.216114"
$ws1.Cells.Item(110, 6).Value = "2021-10-05 14:20:45.216119"
$ws1.Cells.Item(111, 6).Value = "2021-10-05 14:20:45.216122"
$ws1.Cells.Item(112, 6).Value = "2021-10-05 14:20:45.216124"
$ws1.Cells.Item(113, 6).Value = "2021-10-05 14:20:45.216127"
$ws1.Cells.Item(114, 6).Value = "2021-10-05 14:20:45.216129"
$ws1.Cells.Item(115, 6).Value = "2021-10-05 14:20:45.216132"
$ws1.Cells.Item(116, 6).Value = "2021-10-05 14:20:45.216134"
$ws1.Cells.Item(117, 6).Value = "2021-10-05 14:20:45.216137"
$ws1.Cells.Item(118, 6).Value = "2021-10-05 14:20:45.216140"
$ws1.Cells.Item(119, 6).Value = "2021-10-05 14:20:45.216142"
$ws1.Cells.Item(120, 6).Value = "2021-10-05 14:20:45.216145"
$ws1.Cells.Item(121, 6).Value = "2021-10-05 14:20:45.216147"
$ws1.Cells.Item(122, 6).Value = "2021-10-05 14:20:45.216150"
$ws1.Cells.Item(123, 6).Value = "2021-10-05 14:20:45.216152"
$ws1.Cells.Item(124, 6).Value = "2021-10-05 14:20:45.216154"
$ws1.Cells.Item(125, 6).Value = "2021-10-05 14:20:45.216157"
$ws1.Cells.Item(126, 6).Value = "2021-10-05 14:20:45.216160"
$ws1.Cells.Item(127, 6).Value = "2021-10-05 14:20:45.216162"
$ws1.Cells.Item(128, 6).Value = "2021-10-05 14:20:45.216165"
$ws1.Cells.Item(129, 6).Value = "2021-10-05 14:20:45.216167"
$ws1.Cells.Item(130, 6).Value = "2021-10-05 14:20:45.216172"
$ws1.Cells.Item(131, 6).Value = "2021-10-05 14:20:45.216175"
$ws1.Cells.Item(132, 6).Value = "2021-10-05 14:20:45.216177"
$ws1.Cells.Item(133, 6).Value = "2021-10-05 14:20:45.216180"
$ws1.Cells.Item(134, 6).Value = "2021-10-05 14:20:45.216182"
$ws1.Cells.Item(135, 6).Value = "2021-10-05 14:20:45.216185"
$ws1.Cells.Item(136, 6).Value = "2021-10-05 14:20:45.216187"
$ws1.Cells.Item(137, 6).Value = "2021-10-05 14:20:45.216190"
$ws1.Cells.Item(138, 6).Value = "2021-10-05 14:20:45.216192"
$ws1.Cells.Item(139, 6).Value = "2021-10-05 14:20:45.216195"
$ws1.Cells.Item(140, 6).Value = "2021-10-05 14:20:45.216197"
$ws1.Cells.Item(141, 6).Value = "2021-10-05 14:20:45.216200"
$ws1.Cells.Item(142, 6).Value = "2021-10-05 14:20:45.216202"
$ws1.Cells.Item(143, 6).Value = "2021-10-05 14:20:45.216205"
$ws1.Cells.Item(144, 6).Value = "2021-10-05 14:20:45.216207"
$ws1.Cells.Item(145, 6).Value = "2021-10-05 14:20:45.216210"
$ws1.Cells.Item(146, 6).Value = "2021-10-05 14:20:45.216212"
$ws1.Cells.Item(147, 6).Value = "2021-10-05 14:20:45.216215"
$ws1.Cells.Item(148, 6).Value = "2021-10-05 14:20:45.216217"
$ws1.Cells.Item(149, 6).Value = "2021-10-05 14:20:45.216220"
$ws1.Cells.Item(150, 6).Value = "2021-10-05 14:20:45.216223"
$ws1.Cells.Item(151, 6).Value = "2021-10-05 14:20:45.216226"
$ws1.Cells.Item(152, 6).Value = "2021-10-05 14:20:45.216228"
$ws1.Cells.Item(153, 6).Value = "2021-10-05 14:20:45.216231"
$ws1.Cells.Item(154, 6).Value = "2021-10-05 14:20:45.216233"
$ws1.Cells.Item(155, 6).Value = "2021-10-05 14:20:45.216236"
$ws1.Cells.Item(156, 6).Value = "2021-10-05 14:20:45.216238"
$ws1.Cells.Item(157, 6).Value = "2021-10-05 14:20:45.216240"
$ws1.Cells.Item(158, 6).Value = "2021-10-05 14:20:45.216243"
$ws1.Cells.Item(159, 6).Value = "2021-10-05 14:20:45.216245"
$ws1.Cells.Item(160, 6).Value = "2021-10-05 14:20:45.216248"
$ws1.Cells.Item(161, 6).Value = "2021-10-05 14:20:45.216250"
$ws1.Cells.Item(162, 6).Value = "2021-10-05 14:20:45.216253"
$ws1.Cells.Item(163, 6).Value = "2021-10-05 14:20:45.216255"
$ws1.Cells.Item(164, 6).Value = "2021-10-05 14:20:45.216258"
$ws1.Cells.Item(165, 6).Value = "2021-10-05 14:20:45.216260"
$ws1.Cells.Item(166, 6).Value = "2021-10-05 14:20:45.216263"
$ws1.Cells.Item(167, 6).Value = "2021-10-05 14:20:45.216265"
$ws1.Cells.Item(168, 6).Value = "2021-10-05 14:20:45.216268"
$ws1.Cells.Item(169, 6).Value = "2021-10-05 14:20:45.216270"
$ws1.Cells.Item(170, 6).Value = "2021-10-05 14:20:45.216272"
$ws1.Cells.Item(171, 6).Value = "2021-10-05 14:20:45.216275"
$ws1.Cells.Item(172, 6).Value = "2021-10-05 14:20:45.216278"
$ws1.Cells.Item(173, 6).Value = "2021-10-05 14:20:45.216280"
$ws1.Cells.Item(174, 6).Value = "2021-10-05 14:20:45.216284"
$ws1.Cells.Item(175, 6).Value = "2021-10-05 14:20:45.216287"
$ws1.Cells.Item(176, 6).Value = "2021-10-05 14:20:45.216289"
$ws1.Cells.Item(177, 6).Value = "2021-10-05 14:20:45.216292"
$ws1.Cells.Item(178, 6).Value = "2021-10-05 14:20:45.216294"
$ws1.Cells.Item(179, 6).Value = "2021-10-05 14:20:45.216297"
$ws1.Cells.Item(180, 6).Value = "2021-10-05 14:20:45.216299"
$ws1.Cells.Item(181, 6).Value = "2021-10-05 14:20:45.216302"
$ws1.Cells.Item(182, 6).Value = "2021-10-05 14:20:45.216304"
$ws1.Cells.Item(183, 6).Value = "2021-10-05 14:20:45.216307"
$ws1.Cells.Item(184, 6).Value = "2021-10-05 14:20:45.216309"
$ws1.Cells.Item(185, 6).Value = "2021-10-05 14:20:45.216312"
$ws1.Cells.Item(186, 6).Value = "2021-10-05 14:20:45.216314"
$ws1.Cells.Item(187, 6).Value = "2021-10-05 14:20:45.216317"
$ws1.Cells.Item(188, 6).Value = "2021-10-05 14:20:45.216319"
$ws1.Cells.Item(189, 6).Value = "2021-10-05 14:20:45.216322"
$ws1.Cells.Item(190, 6).Value = "2021-10-05 14:20:45.216324"
$ws1.Cells.Item(191, 6).Value = "2021-10-05 14:20:45.216326"
$ws1.Cells.Item(192, 6).Value = "2021-10-05 14:20:45.216329"
$ws1.Cells.Item(193, 6).Value = "2021-10-05 14:20:45.216331"
$ws1.Cells.Item(194, 6).Value = "2021-10-05 14:20:45.216334"
$ws1.Cells.Item(195, 6).Value = "2021-10-05 14:20:45.216336"
$ws1.Cells.Item(196, 6).Value = "2021-10-05 14:20:45.216339"
$ws1.Cells.Item(197, 6).Value = "2021-10-05 14:20:45.216342"
$ws1.Cells.Item(198, 6).Value = "2021-10-05 14:20:45.216344"
$ws1.Cells.Item(199, 6).Value = "2021-10-05 14:20:45.216347"
$ws1.Cells.Item(200, 6).Value = "2021-10-05 14:20:45.216349"
$ws1.Cells.Item(201, 6).Value = "2021-10-05 14:20:45.216352"
$ws1.Cells.Item(202, 6).Value = "2021-10-05 14:20:45.216354"
$ws1.Cells.Item(203, 6).Value = "2021-10-05 14:20:45.216357"
$ws1.Cells.Item(204, 6).Value = "2021-10-05 14:20:45.216360"
$ws1.Cells.Item(205, 6).Value = "2021-10-05 14:20:45.216362"
$ws1.Cells.Item(206, 6).Value = "2021-10-05 14:20:45.216365"
$ws1.Cells.Item(207, 6).Value = "2021-10-05 14:20:45.216367"
$ws1.Cells.Item(208, 6).Value = "2021-10-05 14:20:45.216370"
$ws1.Cells.Item(209, 6).Value = "2021-10-05 14:20:45.216372"
$ws1.Cells.Item(210, 6).Value = "2021-10-05 14:20:45.216375"
$ws1.Cells.Item(211, 6).Value = "2021-10-05 14:20:45.216377"
$ws1.Cells.Item(212, 6).Value = "2021-10-05 14:20:45.216380"
$ws1.Cells.Item(213, 6).Value = "2021-10-05 14:20:45.216382"
$ws1.Cells.Item(214, 6).Value = "2021-10-05 14:20:45.216385"
$ws1.Cells.Item(215, 6).Value = "2021-10-05 14:20:45.216387"
$ws1.Cells.Item(216, 6).Value = "2021-10-05 14:20:45.216390"
$ws1.Cells.Item(217, 6).Value = "2021-10-05 14:20:45.216393"
$ws1.Cells.Item(218, 6).Value = "2021-10-05 14:20:45.216395"
$ws1.Cells.Item(219, 6).Value = "2021-10-05 14:20:45.216398"
$ws1.Cells.Item(220, 6).Value = "2021-10-05 14:20:45.216400"
$ws1.Cells.Item(221, 6).Value = "2021-10-05 14:20:45.216402"
$ws1.Cells.Item(222, 6).Value = "2021-10-05 14:20:45.216405"
$ws1.Cells.Item(223, 6).Value = "2021-10-05 14:20:45.216408"
$ws1.Cells.Item(224, 6).Value = "2021-10-05 14:20:45.216410"
$ws1.Cells.Item(225, 6).Value = "2021-10-05 14:20:45.216414"
$ws1.Cells.Item(226, 6).Value = "2021-10-05 14:20:45.216416"
$ws1.Cells.Item(227, 6).Value = "2021-10-05 14:20:45.216419"
$ws1.Cells.Item(228, 6).Value = "2021-10-05 14:20:45.216421"
$ws1.Cells.Item(229, 6).Value = "2021-10-05 14:20:45.216424"
$ws1.Cells.Item(230, 6).Value = "2021-10-05 14:20:45.216426"
$ws1.Cells.Item(231, 6).Value = "2021-10-05 14:20:45.216429"
$ws1.Cells.Item(232, 6).Value = "2021-10-05 14:20:45.216431"
$ws1.Cells.Item(233, 6).Value = "2021-10-05 14:20:45.216434"
$ws1.Cells.Item(234, 6).Value = "2021-10-05 14:20:45.216437"
$ws1.Cells.Item(235, 6).Value = "2021-10-05 14:20:45.216440"
$ws1.Cells.Item(236, 6).Value = "2021-10-05 14:20:45.216443"
$ws1.Cells.Item(237, 6).Value = "2021-10-05 14:20:45.216446"
$ws1.Cells.Item(238, 6).Value = "2021-10-05 14:20:45.216448"
$ws1.Cells.Item(239, 6).Value = "2021-10-05 14:20:45.216451"
$ws1.Cells.Item(240, 6).Value = "2021-10-05 14:20:45.216453"
$ws1.Cells.Item(241, 6).Value = "2021-10-05 14:20:45.216455"
$ws1.Cells.Item(242, 6).Value = "2021-10-05 14:20:45.216458"
$ws1.Cells.Item(243, 6).Value = "2021-10-05 14:20:45.216460"
$ws1.Cells.Item(244, 6).Value = "2021-10-05 14:20:45.216463"
$ws1.Cells.Item(245, 6).Value = "2021-10-05 14:20:45.216465"
$ws1.Cells.Item(246, 6).Value = "2021-10-05 14:20:45.216468"
$ws1.Cells.Item(247, 6).Value = "2021-10-05 14:20:45.216471"
$ws1.Cells.Item(248, 6).Value = "2021-10-05 14:20:45.216473"
$ws1.Cells.Item(249, 6).Value = "2021-10-05 14:20:45.216476"
$ws1.Cells.Item(250, 6).Value = "2021-10-05 14:20:45.216478"
$ws1.Cells.Item(251, 6).Value = "2021-10-05 14:20:45.216481"
$ws1.Cells.Item(252, 6).Value = "2021-10-05 14:20:45.216484"
$ws1.Cells.Item(253, 6).Value = "2021-10-05 14:20:45.216486"
$ws1.Cells.Item(254, 6).Value = "2021-10-05 14:20:45.216489"
$ws1.Cells.Item(255, 6).Value = "2021-10-05 14:20:45.216491"
$ws1.Cells.Item(256, 6).Value = "2021-10-05 14:20:45.216494"
$ws1.Cells.Item(257, 6).Value = "2021-10-05 14:20:45.216496"
$ws1.Cells.Item(258, 6).Value = "2021-10-05 14:20:45.216499"
$ws1.Cells.Item(259, 6).Value = "2021-10-05 14:20:45.216501"
$ws1.Cells.Item(260, 6).Value = "2021-10-05 14:20:45.216504"
$ws1.Cells.Item(261, 6).Value = "2021-10-05 14:20:45.216506"
$ws1.Cells.Item(262, 6).Value = "2021-10-05 14:20:45.216509"
$ws1.Cells.Item(263, 6).Value = "2021-10-05 14:20:45.216511"
$ws1.Cells.Item(264, 6).Value = "2021-10-05 14:20:45.216514"
$ws1.Cells.Item(265, 6).Value = "2021-10-05 14:20:45.216517"
$ws1.Cells.Item(266, 6).Value = "2021-10-05 14:20:45.216520"
$ws1.Cells.Item(267, 6).Value = "2021-10-05 14:20:45.216522"
$ws1.Cells.Item(268, 6).Value = "2021-10-05 14:20:45.216525"
$ws1.Cells.Item(269, 6).Value = "2021-10-05 14:20:45.216527"
$ws1.Cells.Item(270, 6).Value = "2021-10-05 14:20:45.216530"
$ws1.Cells.Item(271, 6).Value = "2021-10-05 14:20:45.216533"
$ws1.Cells.Item(272, 6).Value = "2021-10-05 14:20:45.216535"
$ws1.Cells.Item(273, 6).Value = "2021-10-05 14:20:45.216538"
$ws1.Cells.Item(274, 6).Value = "2021-10-05 14:20:45.216541"
$ws1.Cells.Item(275, 6).Value = "2021-10-05 14:20:45.216543"
$ws1.Cells.Item(276, 6).Value = "2021-10-05 14:20:45.216546"
$ws1.Cells.Item(277, 6).Value = "2021-10-05 14:20:45.216548"
$ws1.Cells.Item(278, 6).Value = "2021-10-05 14:20:45.216551"
$ws1.Cells.Item(279, 6).Value = "2021-10-05 14:20:45.216554"
$ws1.Cells.Item(280, 6).Value = "2021-10-05 14:20:45.216556"
$ws1.Cells.Item(281, 6).Value = "2021-10-05 14:20:45.216559"
$ws1.Cells.Item(282, 6).Value = "2021-10-05 14:20:45.216562"
$ws1.Cells.Item(283, 6).Value = "2021-10-05 14:20:45.216564"
$ws1.Cells.Item(284, 6).Value = "2021-10-05 14:20:45.216567"

# --- Add the new "metadata" sheet, placed right after "data" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadata"

$ws2.Cells.Item(1, 2).Value = "data_name"
$ws2.Cells.Item(1, 3).Value = "data_id"
$ws2.Cells.Item(1, 4).Value = "data_version"
$ws2.Cells.Item(1, 5).Value = "data_version_created"
$ws2.Cells.Item(1, 6).Value = "panel_query_time"
$ws2.Cells.Item(1, 7).Value = "panel_get_request"

$ws2.Cells.Item(2, 1).Value = 0
$ws2.Cells.Item(2, 2).Value = "Hereditary neuropathy"
$ws2.Cells.Item(2, 3).Value = 85

# data_version ("1.415") must stay text, not become the number 1.415 --
# force Text format just long enough to enter the value as a string, ...
$ws2.Range("D2").NumberFormat = "@"
$ws2.Cells.Item(2, 4).Value = "1.415"
# ... then strip the format back to plain/default so no extra style lingers
$ws1.Range("A1").Copy()
$ws2.Range("D2").PasteSpecial(-4122)

$ws2.Cells.Item(2, 5).Value = "2021-09-09T09:38:04.958233Z"
$ws2.Cells.Item(2, 6).Value = "2021-10-05 14:20:45.212600"
$ws2.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/85/?format=json"

# --- Copy the bold/border/alignment header style from the "data" sheet onto the new headers ---
$ws1.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)
$ws2.Range("A2").PasteSpecial(-4122)

# Keep "data" as the active sheet (matches original activeTab=0)
$ws1.Activate()
